$d = $word.ActiveDocument

# The "ARDUINO CODE (Receiver)" section originally used PIN 11 for TX / PIN 12 for RX
# (and for the hc-05 rx/tx wiring + the "only read from" pin). The fix swaps the two
# pin numbers throughout that section. There are exactly three "PIN 11" and three
# "PIN 12" occurrences in the whole document, all inside this section, so a global
# swap (via a temporary placeholder, to avoid a second pass re-swapping the
# just-written values) reproduces the intended edit without disturbing any other
# run/proofErr structure in the surrounding text.

$placeholder = "@@PIN_SWAP_TMP@@"

# PIN 11 -> placeholder
while ($d.Content.Find.Execute("PIN 11", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "PIN $placeholder", 2)) {}

# PIN 12 -> PIN 11
while ($d.Content.Find.Execute("PIN 12", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "PIN 11", 2)) {}

# placeholder -> 12
while ($d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "12", 2)) {}
